$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new column D values (final parent measurements)
$ws.Range("D2").Value = "Not stressful"
$ws.Range("D3").Value = "Moderately stressful"
$ws.Range("D4").Value = "A little stressful"
$ws.Range("D5").Value = "A little stressful"
$ws.Range("D6").Value = "Moderately stressful"
$ws.Range("D7").Value = "Very stressful"

# Adjust column D width to fit the new content
$ws.Columns.Item(4).ColumnWidth = 17

# Update the active selection to D8, matching where the user left off entering data
$ws.Range("D8").Select()
